# Updated symbol list on Fri Dec 23 11:21:21 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as literal TEXT
# (the source feed writes them with fixed decimal places, e.g. "0.8080",
# "22.00", "0.00005623"). A plain Range.Value assignment of a numeric-
# looking string makes Excel auto-convert it to a real number, which would
# silently drop meaningful trailing zeros / flip to scientific notation.
# Prefixing with a leading apostrophe forces Excel to keep it as text,
# exactly like typing '0.8080 into the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

# ── Price (column D) refreshes ───────────────────────────────────────────
Set-TextValue "D2"  "246.04"
Set-TextValue "D3"  "22.01"
Set-TextValue "D4"  "5.429"
Set-TextValue "D5"  "0.05835"
Set-TextValue "D6"  "3.382"
Set-TextValue "D7"  "6.352"
Set-TextValue "D8"  "0.8080"
Set-TextValue "D9"  "0.9688"
Set-TextValue "D10" "0.1430"
Set-TextValue "D11" "0.07466"
Set-TextValue "D12" "0.03286"
Set-TextValue "D13" "0.03048"
Set-TextValue "D14" "4.170"
Set-TextValue "D15" "0.09398"
Set-TextValue "D16" "0.001589"
Set-TextValue "D18" "0.0005891"
Set-TextValue "D19" "0.006153"
Set-TextValue "D20" "0.004108"
Set-TextValue "D21" "0.0009989"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "3.699"
Set-TextValue "D24" "2.221"
Set-TextValue "D25" "0.3208"
Set-TextValue "D27" "0.0003401"

# ── Rows 41-43: symbol list reshuffled (KickToken/BKEXToken/CEJI rotate) ─
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1078"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002551"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003035"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# ── Remaining Price / Volume(1h) refreshes ───────────────────────────────
Set-TextValue "D44" "0.006687"
Set-TextValue "D45" "0.00005623"

Set-TextValue "D47" "0.4201"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

Set-TextValue "D48" "0.1458"
$ws.Range("E48").Value = "47BOLOBOLO"
